$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "17:09 03-Dec-23"
$ws.Range("C2").Value = "Ẩn danh"
$ws.Range("D2").Value = "A duy"
